$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update score values (GitHub commit stats + basic-option scores) ---
$ws.Range("C9").Value = 28
$ws.Range("C11").Value = 15
$ws.Range("C12").Value = 30
$ws.Range("C16").Value = 8

# --- Move the active cell / selection to F12 ---
$ws.Range("F12").Select() | Out-Null
